# New crime data collected: update header volume/date strings plus the
# Week to Date / 28 Day / Year to Date / historical comparison figures in
# the crime-complaints table (rows 14-29) for the 113th Precinct report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Header: bump the volume/report-number and advance the reporting week.
# ----------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# Helper: xlPasteFormats
$xlPasteFormats = -4122

# ----------------------------------------------------------------------
# Row 14 - Murder
# ----------------------------------------------------------------------
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F14").Value = 1
$ws.Range("J14").Copy() | Out-Null
$ws.Range("F14").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("I14").Value = 1
$ws.Range("J14").Copy() | Out-Null
$ws.Range("I14").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("K14").Value = -66.666666666666
$ws.Range("L14").Value = -66.666666666666
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -87.5

# ----------------------------------------------------------------------
# Row 15 - Rape
# ----------------------------------------------------------------------
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -33.333333333333

# ----------------------------------------------------------------------
# Row 16 - Robbery
# ----------------------------------------------------------------------
$ws.Range("C16").Value = 1

$ws.Range("D16").Value = "'0"
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E16").Value = "'***.*"
$ws.Range("E22").Copy() | Out-Null
$ws.Range("E16").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F16").Value = 11
$ws.Range("H16").Value = 120
$ws.Range("I16").Value = 36
$ws.Range("K16").Value = 44
$ws.Range("L16").Value = -5.263157894736
$ws.Range("M16").Value = -56.626506024096
$ws.Range("N16").Value = -87.878787878787

# ----------------------------------------------------------------------
# Row 17 - Fel. Assault
# ----------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = -19.191919191919
$ws.Range("L17").Value = -13.978494623655
$ws.Range("M17").Value = 5.263157894736
$ws.Range("N17").Value = -52.941176470588

# ----------------------------------------------------------------------
# Row 18 - Burglary
# ----------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -47.058823529411
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = -45.652173913043
$ws.Range("L18").Value = -39.024390243902
$ws.Range("M18").Value = -73.684210526315
$ws.Range("N18").Value = -89.626556016597

# ----------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ----------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 12.903225806451
$ws.Range("I19").Value = 110
$ws.Range("J19").Value = 104
$ws.Range("K19").Value = 5.769230769230
$ws.Range("L19").Value = 37.5
$ws.Range("M19").Value = -29.936305732484
$ws.Range("N19").Value = -85.215053763440

# ----------------------------------------------------------------------
# Row 20 - G.L.A.
# ----------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -37.5
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -44.827586206896
$ws.Range("I20").Value = 51
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -40.697674418604
$ws.Range("L20").Value = 59.375
$ws.Range("M20").Value = -21.538461538461
$ws.Range("N20").Value = -88.221709006928

# ----------------------------------------------------------------------
# Row 21 - TOTAL
# ----------------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -40.540540540540
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -23.728813559322
$ws.Range("I21").Value = 311
$ws.Range("J21").Value = 375
$ws.Range("K21").Value = -17.066666666666
$ws.Range("L21").Value = 5.423728813559
$ws.Range("M21").Value = -35.343035343035
$ws.Range("N21").Value = -83.785192909280

# ----------------------------------------------------------------------
# Row 23 - Housing
# ----------------------------------------------------------------------
$ws.Range("G23").Value = "'0"
$ws.Range("F23").Copy() | Out-Null
$ws.Range("G23").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("H23").Value = "'***.*"
$ws.Range("E23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial($xlPasteFormats) | Out-Null

# ----------------------------------------------------------------------
# Row 24 - Petit Larceny
# ----------------------------------------------------------------------
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = -18.947368421052
$ws.Range("I24").Value = 317
$ws.Range("J24").Value = 360
$ws.Range("K24").Value = -11.944444444444
$ws.Range("L24").Value = 18.726591760299
$ws.Range("M24").Value = 26.294820717131

# ----------------------------------------------------------------------
# Row 25 - Misd. Assault
# ----------------------------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 150
$ws.Range("I25").Value = 164
$ws.Range("J25").Value = 108
$ws.Range("K25").Value = 51.851851851851
$ws.Range("L25").Value = 56.190476190476
$ws.Range("M25").Value = -19.211822660098

# ----------------------------------------------------------------------
# Row 26 - UCR Rape*
# ----------------------------------------------------------------------
$ws.Range("D26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = -14.285714285714

# ----------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ----------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666

# ----------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ----------------------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("F28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F28").Value = 2
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = 133.333333333333
$ws.Range("N28").Value = -84.090909090909

# ----------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ----------------------------------------------------------------------
$ws.Range("C29").Value = 1
$ws.Range("F29").Copy() | Out-Null
$ws.Range("C29").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 66.666666666666
$ws.Range("L29").Value = -28.571428571428
$ws.Range("M29").Value = -58.333333333333
$ws.Range("N29").Value = -86.111111111111
